# Update Total Deaths (column B) values in Sheet1 with newer data, per upstream dataset refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> updated "Total Deaths" value (column B)
$updates = @{
    2 = 7325
    3 = 3132
    4 = 6137
    5 = 133
    6 = 1737
    7 = 117
    8 = 116760
    9 = 7788
    10 = 2104
    11 = 13143
    12 = 8070
    13 = 708
    15 = 28022
    16 = 251
    17 = 5261
    18 = 27504
    19 = 586
    22 = 19305
    23 = 12918
    24 = 2421
    25 = 616457
    26 = 98
    27 = 29486
    28 = 290
    29 = 19177
    31 = 351
    32 = 2982
    33 = 1823
    34 = 29963
    36 = 181
    37 = 38600
    39 = 129056
    40 = 151
    41 = 359
    42 = 1118
    43 = 7332
    44 = 706
    45 = 11574
    46 = 8313
    47 = 608
    48 = 34377
    49 = 3018
    50 = 189
    51 = 42
    52 = 4218
    53 = 33498
    54 = 21015
    55 = 3793
    56 = 175
    57 = 63
    58 = 1848
    60 = 6829
    61 = 697
    62 = 1421
    63 = 121368
    64 = 282
    66 = 12695
    67 = 105642
    68 = 1239
    69 = 19157
    71 = 16036
    72 = 388
    73 = 149
    74 = 1016
    75 = 758
    76 = 10415
    77 = 36429
    78 = 35
    79 = 475434
    80 = 143929
    81 = 130661
    82 = 23965
    83 = 5788
    84 = 8210
    85 = 134765
    86 = 2424
    87 = 18369
    88 = 11941
    89 = 18005
    90 = 5346
    92 = 4253
    93 = 2987
    94 = 2466
    95 = 2768
    96 = 236
    97 = 4366
    98 = 8836
    99 = 664
    101 = 5540
    102 = 67
    103 = 6956
    104 = 890
    105 = 980
    106 = 2307
    107 = 30862
    108 = 257
    109 = 628
    110 = 470
    112 = 849
    114 = 296385
    115 = 9354
    117 = 2033
    118 = 2349
    119 = 14796
    120 = 1942
    121 = 3577
    122 = 11553
    123 = 20530
    124 = 46
    125 = 214
    126 = 267
    127 = 2981
    128 = 1136
    130 = 28830
    132 = 7386
    133 = 573
    134 = 16500
    135 = 201733
    136 = 50096
    137 = 88414
    138 = 18645
    139 = 612
    140 = 57531
    141 = 282713
    142 = 1344
    144 = 283
    145 = 77
    147 = 94
    148 = 57
    149 = 8852
    150 = 1886
    151 = 12148
    152 = 129
    154 = 789
    155 = 15304
    156 = 5394
    158 = 1333
    159 = 90116
    161 = 88381
    162 = 14595
    163 = 3208
    164 = 1176
    165 = 15191
    166 = 11748
    167 = 2809
    170 = 734
    171 = 21171
    175 = 2383
    176 = 25427
    177 = 78969
    178 = 797179
    179 = 3267
    180 = 96427
    181 = 2151
    182 = 146844
    183 = 6142
    184 = 1434
    186 = 5234
    187 = 27611
    188 = 4830
    189 = 1968
    190 = 3669
    191 = 4735
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
